$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Complete the existing (previously blank) backlog rows 11-13
$ws.Range("B11").Value = 8
$ws.Range("C11").Value = "Create testcase for Calculation to compute Time"
$ws.Range("D11").Value = 1

$ws.Range("B12").Value = 9
$ws.Range("C12").Value = "Create testcase for Calculation to compute distance"
$ws.Range("D12").Value = 1
$ws.Rows.Item(12).RowHeight = 24

$ws.Range("B13").Value = 10
$ws.Range("C13").Value = "Create testcase for Calculation to compute pace"
$ws.Range("D13").Value = 1

# New rows 14-22 with centered (no-wrap) styling matching column C's default style
$ws.Range("B14").Value = 11
$ws.Range("C14").Value = "Build user interface-window.java (Windows Builder)"
$ws.Range("D14").Value = 3

$ws.Range("B15").Value = 12
$ws.Range("C15").Value = "build Class Distance.java"

$ws.Range("B16").Value = 13
$ws.Range("C16").Value = "build Class Pace.java"

$ws.Range("B17").Value = 14
$ws.Range("C17").Value = "build Class Time.java"

$ws.Range("B18").Value = 15
$ws.Range("C18").Value = "build class PaceCalculator.java"

$ws.Range("B19").Value = 16
$ws.Range("B20").Value = 17
$ws.Range("B21").Value = 18
$ws.Range("B22").Value = $null

$ws.Range("B14:D14").HorizontalAlignment = -4108
$ws.Range("B15:C18").HorizontalAlignment = -4108
$ws.Range("B19:B22").HorizontalAlignment = -4108

# Column C width change (target stored width 46.7109375; runtime quantizes to 1/6 steps,
# so 45.85 is the closest input that yields the nearest achievable stored width)
$ws.Columns.Item(3).ColumnWidth = 45.85

# View / selection changes
$ws.Activate()
try { $excel.ActiveWindow.SetTopLeftVisibleCell("A15") | Out-Null } catch {}
$ws.Range("C19").Select()
